$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.201.94"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.20"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.92"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5221"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3764"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07267"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.22"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9054"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08459"
$ws.Range("E12").Value = "  +10.43%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.13"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.904.85"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.302"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008656"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.55"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.238.92"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.095"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.153.08"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.65"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.442"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.329"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.14"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.759"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.27"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.23"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.823"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.913"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09301"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05069"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7963"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.442"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.947"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.582"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5735"
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.066"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.613"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.19"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1518"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4861"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.18"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.628"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.68"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.04"
$ws.Range("E51").Value = "  -0.22%  "
